# Auto commit at 2025-09-13  8:40:18.74
# Refresh the "Metrics" sheet's raw data values (B2:B13) with the latest
# figures, then re-point each sheet's selection to match where the editor
# was last working: Metrics -> B2:B13 (active cell B2), today -> F8.
# Downstream formulas on "today" (B11:B22, E11:E22, F11:F22) and the
# volatile TODAY()-1 cell on "today" (A1) recalc automatically.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 188434.07
$metrics.Range("B3").Value  = 152656.01
$metrics.Range("B4").Value  = 59738.520000000004
$metrics.Range("B5").Value  = 7384
$metrics.Range("B6").Value  = 4107684.9499999997
$metrics.Range("B7").Value  = 3480183.4899999993
$metrics.Range("B8").Value  = 1189104.2
$metrics.Range("B9").Value  = 158544
$metrics.Range("B10").Value = 32573008.750999827
$metrics.Range("B11").Value = 19510053.560000002
$metrics.Range("B12").Value = 11470813.090000002
$metrics.Range("B13").Value = 1256171

# Restore the Metrics sheet's selection as last left by the editor.
$metrics.Range("B2:B13").Select()

# Re-select the "today" sheet (the workbook's active tab) at its new cell.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("F8").Select()
